# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right before the existing "2022-Q2"
#    sheet (i.e. right after "总计"), seeded as a copy of "2022-Q2" (same
#    layout/formatting), then overwrite its data with the new Q3 numbers.
# 2. Update the "总计" (summary) sheet: the former Q2/Q1/Q4 rows all shift
#    down by one and a new top row for 2022-Q3 is written.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q3" detail sheet.
# ---------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item("2022-Q2")
$sheetQ2.Copy($sheetQ2) | Out-Null
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# Fund code/name (B2/C2) stay the same; only the numeric-looking text
# columns D:G and the rank H change. D:G are stored as text in the
# source workbook, so force text entry (NumberFormat "@"), then strip
# the number-format back off so the cells keep the original "General"
# style (no s= attribute) while remaining text-typed.
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "17.28"
$newSheet.Range("E2").Value = "29.21"
$newSheet.Range("F2").Value = "1.45"
$newSheet.Range("G2").Value = "0.2506"
$newSheet.Range("D2:G2").ClearFormats()
$newSheet.Range("H2").Value = 6

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend the styled "A" column down to the new row 5 (same style as A4).
$total.Range("A4").Copy() | Out-Null
$total.Range("A5").PasteSpecial(-4122) | Out-Null

# New row 5 = what used to be row 4 ("2021-Q4"): A5=3, C5=3, D5=0.16
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.16

# Row 4 becomes the old row 3 data ("2022-Q1")
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.16

# Row 3 becomes the old row 2 data ("2022-Q2")
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.32

# Row 2 becomes the new "2022-Q3" data
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.25

$total.Activate()

Write-Output "2022-Q3 sheet + 总计 summary updated"
